$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 37, shifting rows 37:70 down to 38:71
$ws.Rows("37").Insert()

# Populate the newly inserted row 37 with the new weekly observation.
# Same dimension values as the former row 37 (now row 38), except the
# date (column D) and volume (column J) which are the new data point.
$ws.Range("A37").Value = 9
$ws.Range("B37").Value = 'Vega Central Mapocho de Santiago'
$ws.Range("C37").Value = 'Metropolitana'
$ws.Range("D37").Value = 44484
$ws.Range("E37").Value = 13
$ws.Range("F37").Value = 100112005
$ws.Range("G37").Value = 'Puerro'
$ws.Range("H37").Value = 'Sin especificar'
$ws.Range("I37").Value = 'Primera'
$ws.Range("J37").Value = 160
$ws.Range("K37").Value = 7000
$ws.Range("L37").Value = 8000
$ws.Range("M37").Value = 7500
$ws.Range("N37").Value = '$/paquete 20 unidades'
$ws.Range("O37").Value = 'Provincia de Chacabuco'
$ws.Range("P37").Value = 375
$ws.Range("Q37").Value = 20
$ws.Range("R37").Value = 'Hortaliza'
